$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Item 16 was a mis-labeled duplicate of "Mentor - Student linking" (item 13).
# Rename it to its correct label.
$ws.Range("A19").Value = "16. Form status list"

# Move the active selection to A19 (matches the author's final cursor position).
$ws.Range("A19").Select()
